$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -7.776000000000001
$ws.Range("C3").Value = -12.141
$ws.Range("D5").Value = -7.718000000000001
$ws.Range("D6").Value = -7.605
$ws.Range("D8").Value = -7.903
$ws.Range("C9").Value = -10.614
$ws.Range("A11").Value = -21.3
$ws.Range("B11").Value = 6.504
$ws.Range("A12").Value = -21.465
$ws.Range("C13").Value = -12.546
$ws.Range("C14").Value = -12.569
$ws.Range("A15").Value = -21.628
$ws.Range("D17").Value = -8.014999999999999
$ws.Range("C19").Value = -12.303
$ws.Range("C21").Value = -12.303
$ws.Range("C22").Value = -12.503
$ws.Range("B23").Value = 7.579000000000001
$ws.Range("C24").Value = -12.257
$ws.Range("C26").Value = -11.913
$ws.Range("A27").Value = -21.455
$ws.Range("D27").Value = -7.995
$ws.Range("A28").Value = -20.648
$ws.Range("B28").Value = 6.478
$ws.Range("A31").Value = -21.603
$ws.Range("A32").Value = -20.986
$ws.Range("B32").Value = 7.234999999999999
$ws.Range("D33").Value = -8.111000000000001
$ws.Range("B34").Value = 7.403999999999999
$ws.Range("A36").Value = -20.974
$ws.Range("B36").Value = 6.778999999999999
$ws.Range("B37").Value = 6.956
$ws.Range("A38").Value = -20.574
$ws.Range("C38").Value = -11.941
$ws.Range("C41").Value = -11.842
$ws.Range("B42").Value = 7.761
$ws.Range("A46").Value = -21.547
$ws.Range("B49").Value = 6.474000000000001
$ws.Range("C52").Value = -11.552
$ws.Range("A54").Value = -20.908
$ws.Range("B54").Value = 6.092000000000001
$ws.Range("A55").Value = -22.018
$ws.Range("D55").Value = -7.83
$ws.Range("A56").Value = -21.923
$ws.Range("C56").Value = -12.717
$ws.Range("D59").Value = -7.946
$ws.Range("A67").Value = -21.359
$ws.Range("A69").Value = -21.391
$ws.Range("D70").Value = -7.763
$ws.Range("C71").Value = -11.305
$ws.Range("A72").Value = -21.017
$ws.Range("C72").Value = -12.49
$ws.Range("A73").Value = -19.974
$ws.Range("B78").Value = 7.823
$ws.Range("C78").Value = -11.42
$ws.Range("B80").Value = 7.393000000000001
$ws.Range("D80").Value = -7.588000000000001
$ws.Range("A83").Value = -21.052
$ws.Range("C83").Value = -12.898
$ws.Range("C85").Value = -12.253
$ws.Range("A86").Value = -21.326
$ws.Range("C86").Value = -13.508
$ws.Range("C90").Value = -10.962
$ws.Range("A91").Value = -20.958
$ws.Range("A93").Value = -21.327
$ws.Range("D95").Value = -7.839
$ws.Range("C96").Value = -12.247
$ws.Range("B97").Value = 5.096000000000001
$ws.Range("D97").Value = -7.792000000000002
$ws.Range("D98").Value = -8.063000000000001
$ws.Range("A99").Value = -20.844
$ws.Range("B99").Value = 6.255
$ws.Range("B100").Value = 4.994999999999999
$ws.Range("B101").Value = 5.606
$ws.Range("D102").Value = -7.725
$ws.Range("C103").Value = -12.727
$ws.Range("A104").Value = -21.164
$ws.Range("A105").Value = -20.434
